$d = $word.ActiveDocument

# Update the date/day heading paragraph.
$d.Content.Find.Execute("2024-04-19 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-04-20 Saturday", 2) | Out-Null

# Update the division-problem table. The table has 20 rows x 5 columns, but
# only every 4th row (1, 5, 9, 13, 17) actually carries visible text - the
# rows in between are blank spacer rows. Address each data cell directly by
# (row, column) so there is no ambiguity from duplicate/overlapping values.
$t = $d.Tables.Item(1)

$updates = @(
    @(1, 1, "25÷8="),
    @(1, 2, "37÷9="),
    @(1, 3, "17÷7="),
    @(1, 4, "54÷5="),
    @(1, 5, "90÷6="),

    @(5, 1, "65÷2="),
    @(5, 2, "51÷5="),
    @(5, 3, "52÷9="),
    @(5, 4, "23÷5="),
    @(5, 5, "51÷8="),

    @(9, 1, "24÷8="),
    @(9, 2, "94÷9="),
    @(9, 3, "99÷6="),
    @(9, 4, "33÷9="),
    @(9, 5, "95÷9="),

    @(13, 1, "87÷7="),
    @(13, 2, "45÷6="),
    @(13, 3, "55÷8="),
    @(13, 4, "96÷7="),
    @(13, 5, "75÷2="),

    @(17, 1, "57÷9="),
    @(17, 2, "30÷5="),
    @(17, 3, "27÷3="),
    @(17, 4, "44÷7="),
    @(17, 5, "61÷9=")
)

foreach ($u in $updates) {
    $rowIdx = $u[0]
    $colIdx = $u[1]
    $newVal = $u[2]
    $t.Cell($rowIdx, $colIdx).Range.Text = $newVal
}

Write-Host "Applied date + $($updates.Count) cell updates"
